$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status column ("Handed back: in sync with en-US") for both locale sheets
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value2 = "Handed back: in sync with en-US"

Write-Output "status done"

# ---------------------------------------------------------------------------
# 2. Latest Target File (I) / Latest Handback File (J) for zh-cn
# ---------------------------------------------------------------------------
$dc50655fName = "dc50655f-ca93-4596-8aeb-64a92162e418.md"
$dc50655fUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd84bc6bb6d8bb58950570be2c7546c2071b580f/e2e/dc50655f-ca93-4596-8aeb-64a92162e418.md"
$ffffbbUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd84bc6bb6d8bb58950570be2c7546c2071b580f/e2e/ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md"
$ffffbbName   = "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md"

$wsZhCn.Range("I2").Value2 = $dc50655fName
$wsZhCn.Range("J2").Value2 = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.zh-cn.xlf"
$wsZhCn.Range("K2").Value2 = "2016-09-07 07:28:37"
$wsZhCn.Range("I3").Value2 = $dc50655fName
$wsZhCn.Range("J3").Value2 = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.zh-cn.xlf"
$wsZhCn.Range("K3").Value2 = "2016-09-07 07:28:37"

Write-Output "zh-cn values done"

# ---------------------------------------------------------------------------
# 3. Hyperlinks for zh-cn: keep the existing A2/A3 links untouched and add the
#    two new "Latest Target File" links (I2, I3) pointing at the same file.
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $dc50655fUrl, [Type]::Missing, [Type]::Missing, $dc50655fName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $dc50655fUrl, [Type]::Missing, [Type]::Missing, $dc50655fName) | Out-Null

Write-Output "zh-cn hyperlinks done"

# ---------------------------------------------------------------------------
# 4. Latest Target File (I) / Latest Handback File (J) / Latest Handback
#    DateTime (K) for de-de
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value2 = $dc50655fName
$wsDeDe.Range("J2").Value2 = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.de-de.xlf"
$wsDeDe.Range("K2").Value2 = "2016-09-07 07:28:45"
$wsDeDe.Range("I3").Value2 = $dc50655fName
$wsDeDe.Range("J3").Value2 = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.de-de.xlf"
$wsDeDe.Range("K3").Value2 = "2016-09-07 07:28:45"

Write-Output "de-de values done"

# ---------------------------------------------------------------------------
# 5. Hyperlinks for de-de: keep existing A2/A3 links untouched and add the
#    two new "Latest Target File" links (I2, I3).
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $dc50655fUrl, [Type]::Missing, [Type]::Missing, $dc50655fName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $dc50655fUrl, [Type]::Missing, [Type]::Missing, $dc50655fName) | Out-Null

Write-Output "de-de hyperlinks done"
